# Add a "Torre" (ward/tower) column to the "Relacion de Solicitudes" report,
# right after "Nombre Completo" and before "Cama " — per commit:
# "se añade la celda de torre a la relacion de solicitud de pacientes"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new blank column before column E ("Cama "), shifting Cama..Bebida
# (E:L) one column to the right (F:M), along with their column-width
# definitions, row-1 banner cells and the scattered single cells in rows
# 13-22.
$ws.Columns("E").Insert()

# Give the new header cell the same banner formatting as its neighbours in
# row 1 (style carried by D1/F1, visually a continuous merged band).
$ws.Range("F1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# New header text/value for the inserted column.
$ws.Range("E2").Value2 = "Torre"

# Re-point the AutoFilter so it covers the new last header column (M).
$ws.AutoFilterMode = $false
$ws.Range("B2:M2").AutoFilter()

# Update the (duplicated) hidden _FilterDatabase defined names so they
# reference the widened header range too.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$B`$2:`$M`$2"
    }
}
